$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G19").Value = 3.6
$ws.Range("I19").Value = 2.45
$ws.Range("J19").Value = 4.5
$ws.Range("L19").Value = 3.5
$ws.Range("M19").Value = 1.18
$ws.Range("N19").Value = 4.5
$ws.Range("Q19").Value = 3.6
$ws.Range("R19").Value = 1.29
$ws.Range("X19").Value = 15
$ws.Range("Y19").Value = 15
$ws.Range("AB19").Value = 51
$ws.Range("AD19").Value = 5.5
$ws.Range("AK19").Value = 26
$ws.Range("AO19").Value = 23
$ws.Range("H20").Value = 2.88
$ws.Range("K20").Value = 1.73
$ws.Range("L20").Value = 3.4
$ws.Range("N20").Value = 4.75
$ws.Range("AC20").Value = 4.75
$ws.Range("AR20").Value = 201
$ws.Range("AT20").Value = 1.91
$ws.Range("M21").Value = 1.13
$ws.Range("N21").Value = 6
$ws.Range("G35").Value = 1.6
$ws.Range("H35").Value = 4.1
$ws.Range("I35").Value = 5.25
$ws.Range("K35").Value = 2.38
$ws.Range("L35").Value = 5
$ws.Range("M35").Value = 1.04
$ws.Range("N35").Value = 13
$ws.Range("S35").Value = 1.33
$ws.Range("T35").Value = 3.25
$ws.Range("W35").Value = 8.5
$ws.Range("X35").Value = 8.5
$ws.Range("AG35").Value = 151
$ws.Range("AN35").Value = 3.75
$ws.Range("AT35").Value = 3.25
$ws.Range("G36").Value = 1.62
$ws.Range("Q36").Value = 2.15
$ws.Range("R36").Value = 1.67
$ws.Range("AC36").Value = 8
$ws.Range("AG36").Value = 451
$ws.Range("AH36").Value = 13
$ws.Range("AL36").Value = 51
$ws.Range("AN36").Value = 3.4
$ws.Range("AU36").Value = 9.5
$ws.Range("N37").Value = 9
$ws.Range("Q38").Value = 2.4
$ws.Range("R38").Value = 1.53
$ws.Range("AI41").Value = 11
$ws.Range("AS41").Value = 301
$ws.Range("G86").Value = 2.2
$ws.Range("H86").Value = 3.35
$ws.Range("I86").Value = 2.8
$ws.Range("W86").Value = 7.6
$ws.Range("X86").Value = 9.75
$ws.Range("Z86").Value = 18
$ws.Range("AB86").Value = 19
$ws.Range("AC86").Value = 11.5
$ws.Range("AD86").Value = 5.9
$ws.Range("AF86").Value = 37
$ws.Range("AH86").Value = 9
$ws.Range("AI86").Value = 13.5
$ws.Range("AJ86").Value = 8.75
$ws.Range("AK86").Value = 27
$ws.Range("AL86").Value = 17.5
$ws.Range("AM86").Value = 21
$ws.Range("AN86").Value = 4.25
$ws.Range("AO86").Value = 11.5
$ws.Range("AP86").Value = 18.5
$ws.Range("AR86").Value = 75
$ws.Range("AT86").Value = 2.87
$ws.Range("AU86").Value = 6.7
$ws.Range("AV86").Value = 55
$ws.Range("AW86").Value = 4.85
$ws.Range("AY86").Value = 20
$ws.Range("BA86").Value = 90
$ws.Range("G87").Value = 2.65
$ws.Range("H87").Value = 2.9
$ws.Range("I87").Value = 2.57
$ws.Range("J87").Value = 3.15
$ws.Range("K87").Value = 2.05
$ws.Range("L87").Value = 3.1
$ws.Range("N87").Value = 11
$ws.Range("O87").Value = 1.18
$ws.Range("P87").Value = 4.2
$ws.Range("Q87").Value = 1.55
$ws.Range("R87").Value = 2.15
$ws.Range("S87").Value = 1.29
$ws.Range("T87").Value = 3.3
$ws.Range("W87").Value = 10.5
$ws.Range("X87").Value = 15
$ws.Range("Y87").Value = 8
$ws.Range("Z87").Value = 29
$ws.Range("AA87").Value = 16
$ws.Range("AB87").Value = 16
$ws.Range("AC87").Value = 12.5
$ws.Range("AD87").Value = 5.5
$ws.Range("AE87").Value = 7.9
$ws.Range("AF87").Value = 22
$ws.Range("AG87").Value = 100
$ws.Range("AH87").Value = 10
$ws.Range("AI87").Value = 14
$ws.Range("AK87").Value = 27
$ws.Range("AL87").Value = 15.5
$ws.Range("AM87").Value = 16
$ws.Range("AN87").Value = 5
$ws.Range("AO87").Value = 14.5
$ws.Range("AP87").Value = 16.5
$ws.Range("AQ87").Value = 60
$ws.Range("AS87").Value = 150
$ws.Range("AT87").Value = 2.95
$ws.Range("AU87").Value = 5.6
$ws.Range("AV87").Value = 35
$ws.Range("AW87").Value = 4.9
$ws.Range("AX87").Value = 14
$ws.Range("AZ87").Value = 60
$ws.Range("BA87").Value = 70
$ws.Range("BC87").Value = 450
$ws.Range("G106").Value = 9.5
$ws.Range("I106").Value = 1.36
$ws.Range("L106").Value = 1.91
$ws.Range("Q106").Value = 2.1
$ws.Range("R106").Value = 1.7
$ws.Range("AW106").Value = 3
